# Update the "gpx_191010_pinar" worksheet:
#  - Row 2 (waypoint 1, the northeast corner) previously had "NA" placeholders
#    for longitude/latitude. We now have real coordinates for that point,
#    sourced from Google Maps.
#  - The note explaining why those coordinates are approximate is expanded
#    with additional context about how they were derived/verified.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -121.906856
$ws.Range("B2").Value = 36.593389999999999

$ws.Range("G2").Value = "I thought I got the waypoint for the northeast corner of the property but I didn't; these values are from google maps (I checked point 15 and it was correct to the 4th decimal)"
